$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.89203699999999
$ws.Range("H2").Value = 107.676111
$ws.Range("I2").Value = 0.9301196142645664
$ws.Range("J2").Value = 0.9301196142645662
$ws.Range("M2").Value = 0.011481
$ws.Range("N2").Value = 0.034443
$ws.Range("O2").Value = 0.0845733396193058
$ws.Range("P2").Value = 0.08457333961930581
$ws.Range("Q2").Value = 0.4120764767969999
$ws.Range("R2").Value = 3.708688291173
$ws.Range("S2").Value = 0.07866332202377488
$ws.Range("T2").Value = 0.07866332202377488

$ws.Range("G3").Value = 35.89203699999999
$ws.Range("H3").Value = 107.676111
$ws.Range("I3").Value = 0.9301196142645664
$ws.Range("J3").Value = 0.9301196142645662
$ws.Range("O3").Value = 0.6656623843479286
$ws.Range("P3").Value = 0.6656623843479286
$ws.Range("Q3").Value = 3.243383923504999
$ws.Range("R3").Value = 29.19045531154499
$ws.Range("S3").Value = 0.6191456401601269
$ws.Range("T3").Value = 0.6191456401601267

$ws.Range("G4").Value = 35.89203699999999
$ws.Range("H4").Value = 107.676111
$ws.Range("I4").Value = 0.9301196142645664
$ws.Range("J4").Value = 0.9301196142645662
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.033906
$ws.Range("N4").Value = 0.101718
$ws.Range("O4").Value = 0.2497642760327657
$ws.Range("P4").Value = 0.2497642760327657
$ws.Range("Q4").Value = 1.216955406522
$ws.Range("R4").Value = 10.952598658698
$ws.Range("S4").Value = 0.2323106520806647
$ws.Range("T4").Value = 0.2323106520806646

$ws.Range("I5").Value = 0.02407019339680195
$ws.Range("J5").Value = 0.02407019339680195
$ws.Range("M5").Value = 0.011481
$ws.Range("N5").Value = 0.034443
$ws.Range("O5").Value = 0.0845733396193058
$ws.Range("P5").Value = 0.08457333961930581
$ws.Range("Q5").Value = 0.010663962289
$ws.Range("R5").Value = 0.095975660601
$ws.Range("S5").Value = 0.002035696640850103
$ws.Range("T5").Value = 0.002035696640850103

$ws.Range("I6").Value = 0.02407019339680195
$ws.Range("J6").Value = 0.02407019339680195
$ws.Range("O6").Value = 0.6656623843479286
$ws.Range("P6").Value = 0.6656623843479286
$ws.Range("S6").Value = 0.01602262232823095
$ws.Range("T6").Value = 0.01602262232823095

$ws.Range("I7").Value = 0.02407019339680195
$ws.Range("J7").Value = 0.02407019339680195
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.033906
$ws.Range("N7").Value = 0.101718
$ws.Range("O7").Value = 0.2497642760327657
$ws.Range("P7").Value = 0.2497642760327657
$ws.Range("Q7").Value = 0.031493102114
$ws.Range("R7").Value = 0.283437919026
$ws.Range("S7").Value = 0.006011874427720896
$ws.Range("T7").Value = 0.006011874427720896

$ws.Range("G8").Value = 1.767752333333333
$ws.Range("H8").Value = 5.303257
$ws.Range("I8").Value = 0.04581019233863175
$ws.Range("J8").Value = 0.04581019233863175
$ws.Range("M8").Value = 0.011481
$ws.Range("N8").Value = 0.034443
$ws.Range("O8").Value = 0.0845733396193058
$ws.Range("P8").Value = 0.08457333961930581
$ws.Range("Q8").Value = 0.020295564539
$ws.Range("R8").Value = 0.182660080851
$ws.Range("S8").Value = 0.003874320954680824
$ws.Range("T8").Value = 0.003874320954680824

$ws.Range("G9").Value = 1.767752333333333
$ws.Range("H9").Value = 5.303257
$ws.Range("I9").Value = 0.04581019233863175
$ws.Range("J9").Value = 0.04581019233863175
$ws.Range("O9").Value = 0.6656623843479286
$ws.Range("P9").Value = 0.6656623843479286
$ws.Range("Q9").Value = 0.1597429396016666
$ws.Range("R9").Value = 1.437686456415
$ws.Range("S9").Value = 0.03049412185957082
$ws.Range("T9").Value = 0.03049412185957082

$ws.Range("G10").Value = 1.767752333333333
$ws.Range("H10").Value = 5.303257
$ws.Range("I10").Value = 0.04581019233863175
$ws.Range("J10").Value = 0.04581019233863175
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.033906
$ws.Range("N10").Value = 0.101718
$ws.Range("O10").Value = 0.2497642760327657
$ws.Range("P10").Value = 0.2497642760327657
$ws.Range("Q10").Value = 0.059937410614
$ws.Range("R10").Value = 0.539436695526
$ws.Range("S10").Value = 0.01144174952438011
$ws.Range("T10").Value = 0.01144174952438011
